$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.754.81'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.295.31'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.40%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '98.80'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '270.89'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.26%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.13'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0929'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.89'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.79%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.82'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.639.88'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.848'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.283.08'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.777.83'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.18'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +7.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.08'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.86'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +13.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.10'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.95%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.26'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.05'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '176.53'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.86'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0891'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.43'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.127'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +7.95%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.53'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.83%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.94%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.20'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.70'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.86'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.02%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.79%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '98.37'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.09%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.53'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +12.26%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.441'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.71%  '
